# "ran odds, other fixes"
# Update this week's point spreads ("Week 10" sheet, column D) with the
# latest lines, clearing the leftover explicit cell style on each edited
# cell (Excel drops the manual font override once the cell is re-entered),
# and leave the selection where the user ended up after doing the work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 10")

# New spreads (row -> value), row 3 (Cardinals/Seahawks) is untouched.
$newSpreads = @{
    4  = -2.5
    5  = -1.5
    6  = 10.5
    7  = -1
    8  = -6
    9  = -3.5
    10 = 2
    11 = 4
    12 = -10
    13 = -1.5
    14 = -7
    15 = 7.5
    16 = -4.5
}

foreach ($row in $newSpreads.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $newSpreads[$row]
    $cell.ClearFormats()
}

[void]$ws.Range("S15").Select()
